# Add a new "Identification" data row (row 3) to the VEMS sheet, mirroring
# the layout of the existing row 2 (converter now writes a second record
# after opening/parsing another folder/file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT (not auto-converted by Excel to a
# number/date), then drop the temporary "@" number format so the cell is
# left on the default style - matching how the other text cells on the
# sheet (e.g. row 2) carry no explicit style.
function Set-TextValue {
    param($CellRef, $Text)
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

Set-TextValue "A3" "2"
$ws.Range("B3").Value = "Identification:"
$ws.Range("C3").Value = "29/10/2004"
Set-TextValue "D3" "07.03.22"
$ws.Range("E3").Value = "176 cm"
$ws.Range("F3").Value = "68.0 kg"
Set-TextValue "G3" "22"
Set-TextValue "H3" "3.55"
Set-TextValue "I3" "-0.41"
$ws.Range("J3").Value = 0
Set-TextValue "K3" "63"
$ws.Range("L3").Value = 0
Set-TextValue "M3" "1.45"
$ws.Range("N3").Value = 0
